# Added GEF-based filter sorting columns "bovenkant filter (gef)" (AM) and
# "Filternummer (gef)" (AN). "Filternummer" as originally supplied by GEF
# files comes in a random order, so these two new columns record the GEF
# top-of-filter elevation and the 1/2/3 rank that elevation sorts to
# within each well (group of 3 rows). A couple of previously swapped
# "Zandvanglengte (meters)" values (column Y) are corrected as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to remain plain text
# (keeps parity with other already-existing text cells in this sheet,
# which store numeric-looking strings such as "0.600000" as text) and
# without leaving any extra/unused style definitions or stray cells
# behind. A scratch cell well outside the used range is formatted as
# Text, given the value, then only its *value* (not its number format)
# is pasted onto the real target cell; the scratch column is then
# removed again so it does not affect the sheet's used range/dimension.
function Set-TextValue($row, $col, $text) {
    $helper = $ws.Cells.Item($row, 100)
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)  # xlPasteValues
    $ws.Columns.Item(100).Delete()
}

# --- New header cells (AM1 / AN1), matching the style of existing headers ---
$ws.Range("AL1").Copy()
$ws.Range("AM1:AN1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(1, 39).Value = "bovenkant filter (gef)"
$ws.Cells.Item(1, 40).Value = "Filternummer (gef)"

# --- New data columns: AM = bovenkant filter (gef), AN = Filternummer (gef) ---
$ws.Cells.Item(2, 39).Value = 3.24
$ws.Cells.Item(2, 40).Value = 1
$ws.Cells.Item(3, 39).Value = 9.24
$ws.Cells.Item(3, 40).Value = 2
$ws.Cells.Item(4, 39).Value = 21.24
$ws.Cells.Item(4, 40).Value = 3
$ws.Cells.Item(5, 39).Value = 3.25
$ws.Cells.Item(5, 40).Value = 1
$ws.Cells.Item(6, 39).Value = 9.25
$ws.Cells.Item(6, 40).Value = 2
$ws.Cells.Item(7, 39).Value = 24.25
$ws.Cells.Item(7, 40).Value = 3
$ws.Cells.Item(8, 39).Value = 2.87
$ws.Cells.Item(8, 40).Value = 1
$ws.Cells.Item(9, 39).Value = 8.869999999999999
$ws.Cells.Item(9, 40).Value = 2
$ws.Cells.Item(10, 39).Value = 28.870001
$ws.Cells.Item(10, 40).Value = 3
$ws.Cells.Item(11, 39).Value = -9.720000000000001
$ws.Cells.Item(11, 40).Value = 1
$ws.Cells.Item(12, 39).Value = 0.48
$ws.Cells.Item(12, 40).Value = 2
$ws.Cells.Item(13, 39).Value = 11.23
$ws.Cells.Item(13, 40).Value = 3
$ws.Cells.Item(14, 39).Value = -3.72
$ws.Cells.Item(14, 40).Value = 1
$ws.Cells.Item(15, 39).Value = 3.28
$ws.Cells.Item(15, 40).Value = 2
$ws.Cells.Item(16, 39).Value = 17.030001
$ws.Cells.Item(16, 40).Value = 3
$ws.Cells.Item(17, 39).Value = 1.59
$ws.Cells.Item(17, 40).Value = 1
$ws.Cells.Item(18, 39).Value = 8.390000000000001
$ws.Cells.Item(18, 40).Value = 2
$ws.Cells.Item(19, 39).Value = 28.389999
$ws.Cells.Item(19, 40).Value = 3
$ws.Cells.Item(20, 39).Value = 1.97
$ws.Cells.Item(20, 40).Value = 1
$ws.Cells.Item(21, 39).Value = 7.97
$ws.Cells.Item(21, 40).Value = 2
$ws.Cells.Item(22, 39).Value = 23.969999
$ws.Cells.Item(22, 40).Value = 3
$ws.Cells.Item(23, 39).Value = -0.88
$ws.Cells.Item(23, 40).Value = 1
$ws.Cells.Item(24, 39).Value = 6.12
$ws.Cells.Item(24, 40).Value = 2
$ws.Cells.Item(25, 39).Value = 21.120001
$ws.Cells.Item(25, 40).Value = 3
$ws.Cells.Item(26, 39).Value = -0.3
$ws.Cells.Item(26, 40).Value = 1
$ws.Cells.Item(27, 39).Value = 6.5
$ws.Cells.Item(27, 40).Value = 2
$ws.Cells.Item(28, 39).Value = 21
$ws.Cells.Item(28, 40).Value = 3

# --- Fix swapped Zandvanglengte (meters) values (column Y) ---
Set-TextValue 2 25 "0.630000"
Set-TextValue 4 25 "0.600000"
Set-TextValue 20 25 "0.650000"
Set-TextValue 22 25 "0.600000"
